$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BWTracker")

$task = "import data to database"
$info = "import data from 2 excel file ENGRVW & AR to odc database"
$key = "AscenX"
$parentKey = "SETeam"
$moreInfo = "odcDatabase"

$rows = @(
    @{ row=6;  F=7;  G=23; H=51; I=9;  J=52; K=2;  L="02:28:10" },
    @{ row=7;  F=9;  G=59; H=1;  I=10; J=29; K=11; L="00:30:10" },
    @{ row=8;  F=10; G=40; H=49; I=10; J=56; K=32; L="00:15:42" },
    @{ row=9;  F=11; G=12; H=53; I=12; J=36; K=31; L="01:23:38" },
    @{ row=10; F=13; G=14; H=48; I=13; J=58; K=57; L="00:44:08" },
    @{ row=11; F=14; G=14; H=27; I=15; J=7;  K=32; L="00:53:04" },
    @{ row=12; F=15; G=17; H=45; I=16; J=37; K=38; L="01:19:52" }
)

foreach ($r in $rows) {
    $n = $r.row
    $ws.Cells.Item($n, 1).Value = "Tuesday"
    $ws.Cells.Item($n, 2).Value = 11
    $ws.Cells.Item($n, 3).Value = "July"
    $ws.Cells.Item($n, 4).Value = 2017
    $ws.Cells.Item($n, 5).Value = 28
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I
    $ws.Cells.Item($n, 10).Value = $r.J
    $ws.Cells.Item($n, 11).Value = $r.K
    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $task
    $ws.Cells.Item($n, 14).Value = $info
    $ws.Cells.Item($n, 15).Value = $key
    $ws.Cells.Item($n, 16).Value = $parentKey
    $ws.Cells.Item($n, 17).Value = $moreInfo
    $ws.Rows.Item($n).RowHeight = 30
}
